$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Colab에서 구글 드라이브(Google Drive) 압축파일 다운로드 쉽게하기 (코드 3줄!)"
$ws.Range("E4").Value = "https://teddylee777.github.io/colab/gdrive-dataset"

$ws.Range("D9").Value = "Why SIAI – 3. 박사과정 중 필요성을 느꼈지만 엄두를 못냈던 지식들이라는 확신이 들었습니다"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/why-siai-3-necessary-knowledge-in-phd/#utm_source=rss&utm_medium=rss&utm_campaign=why-siai-3-necessary-knowledge-in-phd"

$ws.Range("D29").Value = "[만화] 인턴일기 66~71"
$ws.Range("E29").Value = "https://blog.promedius.ai/intern-life-10/"

$ws.Range("D37").Value = "[Paper Review]ON CONCEPT-BASED EXPLANATIONS  IN DEEP NEURAL NETWORKS"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1891&mod=document&pageid=1"
